# The "ZENITH (ALN-AGT01-008)" trial was removed from the Query1 table on Sheet1.
# Locate it by value (rather than assuming a fixed row number) and delete the
# whole sheet row: this shifts every row below it up by one, shrinks the
# table/autofilter range to match, and lets the now-unused shared string drop
# out on save.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$target = $ws.Cells.Find("ZENITH (ALN-AGT01-008)")
if ($target -ne $null) {
    $ws.Rows.Item($target.Row).Delete()
}

# The ExternalData_1 workbook-level defined name (used by the query table) still
# points at the old, one-row-taller range after the delete, so bring it in line
# with the table's new extent.
$lo = $ws.ListObjects.Item(1)
$wb.Names.Item("ExternalData_1").RefersTo = "=Sheet1!" + $lo.Range.Address()
